$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) values that look like plain numbers must be forced back to
# text (matching the source data which stores them as literal strings), since
# plain assignment would otherwise let Excel auto-convert them to numbers.

$ws.Range("D2").Value = "66.913.66"
$ws.Range("E2").Value = "  +0.84%  "
$ws.Range("D3").Value = "2.523.29"
$ws.Range("E3").Value = "  -1.93%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "588.15"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "171.69"
$ws.Range("E6").Value = "  +4.15%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.528"
$ws.Range("E8").Value = "  +0.65%  "
$ws.Range("D9").Value = "2.521.51"
$ws.Range("E9").Value = "  -1.97%  "
$ws.Range("E10").Value = "  +0.13%  "
$ws.Range("E11").Value = "  +1.82%  "
$ws.Range("E12").Value = "  -0.68%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.341"
$ws.Range("E13").Value = "  -4.72%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.53"
$ws.Range("E14").Value = "  -0.81%  "
$ws.Range("D15").Value = "3.005.00"
$ws.Range("E15").Value = "  -1.18%  "
$ws.Range("E16").Value = "  -1.18%  "
$ws.Range("D17").Value = "66.825.94"
$ws.Range("E17").Value = "  +0.93%  "
$ws.Range("D18").Value = "2.530.56"
$ws.Range("E18").Value = "  -2.33%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.96"
$ws.Range("E19").Value = "  +3.03%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.28"
$ws.Range("E20").Value = "  -1.27%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "353.70"
$ws.Range("E21").Value = "  +0.92%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.16"
$ws.Range("E22").Value = "  -1.37%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.58"
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("E24").Value = "  +4.93%  "
$ws.Range("E25").Value = "  -0.02%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "69.61"
$ws.Range("E26").Value = "  +1.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.98"
$ws.Range("E27").Value = "  -1.16%  "
$ws.Range("D28").Value = "2.669.08"
$ws.Range("E28").Value = "  -1.41%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.03%  "
$ws.Range("D30").Value = "0.0₃0975"
$ws.Range("E30").Value = "  -0.72%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "531.45"
$ws.Range("E31").Value = "  -0.19%  "
$ws.Range("E32").Value = "  +1.73%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.33"
$ws.Range("E33").Value = "  +0.46%  "
$ws.Range("E34").Value = "  -0.38%  "
$ws.Range("E35").Value = "  -0.82%  "
$ws.Range("E36").Value = "  -0.09%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "157.58"
$ws.Range("E37").Value = "  +0.48%  "
$ws.Range("E38").Value = "  -0.50%  "
$ws.Range("E39").Value = "  -0.92%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.44"
$ws.Range("E40").Value = "  +1.20%  "
$ws.Range("E41").Value = "  -1.65%  "
$ws.Range("E42").Value = "  +0.09%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.07"
$ws.Range("E43").Value = "  -0.25%  "
$ws.Range("E44").Value = "  -0.04%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.48"
$ws.Range("E45").Value = "  +4.23%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "39.68"
$ws.Range("E46").Value = "  -0.65%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "148.56"
$ws.Range("E47").Value = "  +0.07%  "
$ws.Range("E48").Value = "  -1.57%  "
$ws.Range("D49").Value = "0.0₆0275"
$ws.Range("E49").Value = "  -3.65%  "
$ws.Range("E50").Value = "  -1.33%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.69"
$ws.Range("E51").Value = "  +0.25%  "
